$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The scheduled task re-wrote the previous last row's timestamp with a
# slightly more precise value when it re-ran.
$ws.Cells.Item(6, 1).Value = 45864.41694890046

# Append the new reading captured by the scheduled task (row 7), copying
# the date/time number format already used for column A.
$ws.Cells.Item(7, 1).NumberFormat = $ws.Cells.Item(6, 1).NumberFormat
$ws.Cells.Item(7, 1).Value = 45864.45857037041
$ws.Cells.Item(7, 2).Value = 2025
$ws.Cells.Item(7, 3).Value = 30
$ws.Cells.Item(7, 4).Value = 16.83
$ws.Cells.Item(7, 5).Value = 80.23
$ws.Cells.Item(7, 6).Value = 571.29
$ws.Cells.Item(7, 7).Value = 6.31
$ws.Cells.Item(7, 8).Value = "SE"
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = "11:00:20"
